$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font + fill + border) from C1 onto the new
# header cells, then overwrite their values with the new column headers.
$ws.Range("C1").Copy($ws.Range("D1:F1"))
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "TestStartTime"
$ws.Range("F1").Value = "TestEndTime"

# Give the two new columns (TestStartTime / TestEndTime) an explicit,
# best-fit-like custom width, matching the widths used for those headers.
$ws.Columns.Item(5).ColumnWidth = 12.6
$ws.Columns.Item(6).ColumnWidth = 11.6

# Move the active selection to C14.
$ws.Range("C14").Select() | Out-Null
